# Add a new "2022-Q4" quarter sheet to the workbook and record it in the
# "总计" (totals) summary sheet.
#
# The new sheet is created by copying the existing "2022-Q3" sheet (so it
# inherits the same headers/formatting) and placing it immediately before
# "2022-Q3" in tab order, then renaming it and filling in the new quarter's
# figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q4" worksheet just before "2022-Q3".
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2022-Q3")
$templateSheet.Copy($templateSheet)
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

# Fund-level figures for the new quarter (same fund as every other quarter
# on this workbook: 513080 / 华安法国CAC40ETF（QDII）).
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.64"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "93.56"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "6.63"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0424"
$newSheet.Range("H2").Value = 3

# Restore the originally-selected tab (last sheet, "2020-Q4") now that the
# copy operation made the new sheet the active one.
$wb.Worksheets.Item("2020-Q4").Activate()

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a row for 2022-Q4 at the top
#    of the data (row 2) and push the rest down.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Rows.Item(2).Insert()

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 0.04

# Renumber the running index in column A for the rows that shifted down.
for ($r = 3; $r -le 10; $r++) {
    $totals.Cells.Item($r, 1).Value = $r - 2
}
